{"js": "const body = context.document.body;\n\n// Locate the target sentence via a stable, unique anchor phrase before making\n// any edits (search text reflects the *current* document contents, so we run\n// every search before the text around it has been changed by this script).\nconst requiredResults = body.search(\"required stupidest human ever was infected\", { matchCase: true });\nrequiredResults.load(\"text\");\n\nconst andThereforeResults = body.search(\" \\u2013 and therefore you are first deployed there by the UN\", { matchCase: true });\nandThereforeResults.load(\"text\");\n\nawait context.sync();\n\nif (requiredResults.items.length !== 1) {\n  throw new Error(\"Expected exactly one match for the 'required stupidest...' anchor, found \" + requiredResults.items.length);\n}\nif (andThereforeResults.items.length !== 1) {\n  throw new Error(\"Expected exactly one match for the ' - and therefore...' anchor, found \" + andThereforeResults.items.length);\n}\n\n// 1) \"required stupidest human ever was infected\"\n//      -> \"required that patient zero of the zombie pandemic was the stupidest human ever\"\nrequiredResults.items[0].insertText(\n  \"required that patient zero of the zombie pandemic was the stupidest human ever\",\n  \"Replace\"\n);\n\n// 2) \" \\u2013 and therefore you are first deployed there by the UN\"\n//      -> \"; therefore, as the area with the highest concentration of zombies, you are first deployed there by the UN\"\nandThereforeResults.items[0].insertText(\n  \"; therefore, as the area with the highest concentration of zombies, you are first deployed there by the UN\",\n  \"Replace\"\n);\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# The sentence lives in its own paragraph, split across two runs:\n#   \"...meant that it was required stupidest human ever was infected \u2013 so the\n#    zombie apocalypse had to begin in Florida\" + \" \u2013 and therefore you are\n#    first deployed there by the UN\"\n# It becomes:\n#   \"...meant that it was required that patient zero of the zombie pandemic\n#    was the stupidest human ever \u2013 so the zombie apocalypse had to begin in\n#    Florida; therefore, as the area with the highest concentration of\n#    zombies, you are first deployed there by the UN\"\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$found = $rng.Find.Execute(\n    \"required stupidest human ever was infected \u2013 so the zombie apocalypse had to begin in Florida \u2013 and therefore you are first deployed there by the UN\"\n)\nif (-not $found) {\n    throw \"Could not locate the target sentence to edit\"\n}\n\n$rng.Text = \"required that patient zero of the zombie pandemic was the stupidest human ever \u2013 so the zombie apocalypse had to begin in Florida; therefore, as the area with the highest concentration of zombies, you are first deployed there by the UN\"\n"}
